# 22255 - QN Enhancement: update the zip changeset number shown in the
# runbook step from 50770 to 50775 (fixed unable to view warning logs).
#
# The changeset number lives in its own bold/16pt run immediately after a
# run containing "Changeset ". We locate the whole "50770" token, then
# replace just its last two characters ("70") with "75" so the leading
# "507" portion of the number is left completely untouched.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$found = $find.Execute("50770", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate changeset number '50770' in the document."
}

$token = $find.Parent
$tokenStart = $token.Start
$tokenEnd = $token.End

# Range covering just the trailing "70" of "50770".
$tail = $d.Range($tokenEnd - 2, $tokenEnd)
$tail.Text = "7"
$tail.InsertAfter("5")

Write-Output "Changeset number now reads: $($d.Range($tokenStart, $tokenEnd).Text)"
